# Noi dung ngay hoc 16 - Form Upload PHP
# Adds a new row (STT 8, "Thanh Trung") right after the "Van Hien" (STT 7) row
# in the feedback table, keeping the existing rows untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Sanity: the table should currently have 8 rows (header + 7 students);
# the new row is appended at the end, which lands right after the
# "Van Hien" row (the last row of the table).
$newRow = $t.Rows.Add()

# STT
$t.Cell($t.Rows.Count, 1).Range.Text = "8"

# Ho ten
$t.Cell($t.Rows.Count, 2).Range.Text = "Thành Trung"

# Nhan xet
$nhanXet = "+ Tư duy xử lý bài toán tốt`r" + `
           "+ Đáp ứng yêu cầu > 95%`r"
$t.Cell($t.Rows.Count, 3).Range.Text = $nhanXet

# Can chu y
$canChuY = "+ Bài 1, 2, 3 trong hàm cần tư duy sử dụng từ khóa return sao cho hợp lý, cụ thể trong trường hợp này sẽ return 1 kiểu số (int/float)`r" + `
           "+ Bài 7, 8 trong hàm cần tư duy sử dụng từ khóa return sao cho hợp lý, cụ thể trong trường hợp này sẽ return 1 kiểu string`r" + `
           "+ Bài 9 có thể sử dụng cú pháp viết tắt của thẻ for là <?php for(): ?> <?php endfor?> khi viết lồng với HTML, để giảm bớt độ phức tạp khi HTML quá nhiều`r"
$t.Cell($t.Rows.Count, 4).Range.Text = $canChuY
